$d = $word.ActiveDocument
$vtab = [char]11
$oAcute = [char]0xF3

# --------------------------------------------------------------------
# Change 1: the paragraph "<w:br/>3.2.2 RF2. Registró de usuarios" is
# currently split across two <w:r> runs (one holding only the <w:br/>,
# one holding the text) even though both share identical run
# properties. Collapse it to a single run, i.e.
#   <w:r><w:rPr>...</w:rPr><w:br/></w:r><w:r><w:rPr>...</w:rPr><w:t>..</w:t></w:r>
# becomes
#   <w:r><w:rPr>...</w:rPr><w:br/><w:t>..</w:t></w:r>
# --------------------------------------------------------------------
$needleText = "3.2.2 RF2. Registr" + $oAcute + " de usuarios"
$targetParaRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $pText = $para.Range.Text
    # Paragraph.Range.Text includes the trailing paragraph mark; strip
    # it before comparing against our (mark-less) needle.
    if ($pText.Length -gt 0) {
        $pText = $pText.Substring(0, $pText.Length - 1)
    }
    if ($pText -eq ($vtab + $needleText)) {
        $targetParaRange = $para.Range
        break
    }
}

if ($targetParaRange -ne $null) {
    # Exclude the trailing paragraph mark from the range we rewrite.
    $paraRange = $d.Range($targetParaRange.Start, $targetParaRange.End - 1)

    # A same-text assignment is a no-op for the underlying run layout,
    # so first swap in a placeholder (forcing a real content change),
    # then write the original text back. The resulting run is rebuilt
    # from scratch, merging the leading <w:br/> and the text into one
    # <w:r>, just like Word does when it lays out edited runs.
    $paraRange.Text = ($vtab + "TEMP_PLACEHOLDER_TOKEN")

    $full = $d.Content.Text
    $placeholderIdx = $full.IndexOf($vtab + "TEMP_PLACEHOLDER_TOKEN")
    $placeholderLen = ("TEMP_PLACEHOLDER_TOKEN").Length + 1
    $rebuiltRange = $d.Range($placeholderIdx, $placeholderIdx + $placeholderLen)
    $rebuiltRange.Text = ($vtab + $needleText)
}

# --------------------------------------------------------------------
# Change 2: "... implementado para dos terminales, por lo tanto ..."
# becomes "... implementado para cuatro terminales, por lo tanto ...".
# Word splits the run being edited away from its unedited neighbours,
# so the paragraph's single run becomes three runs (identical
# formatting, just split at the edit boundaries).
# --------------------------------------------------------------------
$full = $d.Content.Text
$needleOld = "dos terminales"
$oldIdx = $full.IndexOf($needleOld)
if ($oldIdx -ge 0) {
    $oldWordRange = $d.Range($oldIdx, $oldIdx + 3)

    # Briefly toggle formatting on just the replaced word so the
    # engine splits it into its own run, then reset the formatting so
    # its rPr ends up identical to its (unchanged) neighbours.
    $oldWordRange.Font.Bold = 1
    $oldWordRange.Text = "cuatro"

    $full2 = $d.Content.Text
    $newIdx = $full2.IndexOf("cuatro")
    $newWordRange = $d.Range($newIdx, $newIdx + ("cuatro").Length)
    $newWordRange.Font.Bold = 0
}
